$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tagData = @(
    @(1, "tag", "korTag"),
    @(2, "personal_care", "개인위생/관리용품"),
    @(3, "clothing", "의류"),
    @(4, "personal_care", "개인위생/관리용품"),
    @(5, "personal_care", "개인위생/관리용품"),
    @(6, "baby", "유아용품"),
    @(7, "clothing", "의류"),
    @(8, "personal_care", "개인위생/관리용품"),
    @(9, "stationery", "문구류"),
    @(10, "stationery", "문구류"),
    @(11, "personal_care", "개인위생/관리용품"),
    @(12, "personal_care", "개인위생/관리용품"),
    @(13, "food", "식품"),
    @(14, "food", "식품"),
    @(15, "clothing", "의류"),
    @(16, "medical", "의약/의료"),
    @(17, "baby,food", "유아용품,식품"),
    @(18, "food", "식품"),
    @(19, "clothing", "의류"),
    @(20, "food", "식품"),
    @(21, "medical", "의약/의료"),
    @(22, "personal_care", "개인위생/관리용품"),
    @(23, "food", "식품"),
    @(24, "electronics", "전자기기"),
    @(25, "clothing", "의류"),
    @(26, "clothing", "의류"),
    @(27, "personal_care", "개인위생/관리용품"),
    @(28, "baby,clothing", "유아용품,의류"),
    @(29, "medical", "의약/의료"),
    @(30, "food", "식품"),
    @(31, "clothing", "의류"),
    @(32, "stationery", "문구류"),
    @(33, "stationery", "문구류"),
    @(34, "clothing", "의류"),
    @(35, "personal_care", "개인위생/관리용품"),
    @(36, "electronics", "전자기기"),
    @(37, "clothing", "의류"),
    @(38, "clothing", "의류"),
    @(39, "clothing", "의류"),
    @(40, "food", "식품"),
    @(41, "food", "식품"),
    @(42, "clothing", "의류"),
    @(43, "electronics", "전자기기"),
    @(44, "personal_care", "개인위생/관리용품"),
    @(45, "food", "식품"),
    @(46, "personal_care", "개인위생/관리용품"),
    @(47, "stationery", "문구류"),
    @(48, "clothing", "의류"),
    @(49, "baby", "유아용품"),
    @(50, "food", "식품"),
    @(51, "baby", "유아용품"),
    @(52, "personal_care", "개인위생/관리용품"),
    @(53, "food", "식품"),
    @(54, "medical", "의약/의료"),
    @(55, "clothing", "의류"),
    @(56, "electronics", "전자기기"),
    @(57, "electronics", "전자기기"),
    @(58, "electronics", "전자기기"),
    @(59, "personal_care", "개인위생/관리용품"),
    @(60, "electronics", "전자기기"),
    @(61, "personal_care", "개인위생/관리용품"),
    @(62, "food", "식품"),
    @(63, "electronics", "전자기기"),
    @(64, "clothing", "의류"),
    @(65, "personal_care", "개인위생/관리용품"),
    @(66, "personal_care", "개인위생/관리용품"),
    @(67, "personal_care", "개인위생/관리용품"),
    @(68, "food", "식품"),
    @(69, "medical", "의약/의료"),
    @(70, "clothing", "의류"),
    @(71, "food", "식품"),
    @(72, "clothing", "의류"),
    @(73, "food", "식품"),
    @(74, "clothing", "의류"),
    @(75, "electronics", "전자기기"),
)

foreach ($entry in $tagData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 6).Value = $entry[1]
    $ws.Cells.Item($r, 7).Value = $entry[2]
}

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$ws.Range("F4").Select()
